$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 15000
$ws.Range("I13").Value = 15000
$ws.Range("K13").Value = 15000
$ws.Range("M13").Value = -14831
$ws.Range("H32").Value = 1713
$ws.Range("I32").Value = 1725
$ws.Range("J32").Value = 1701
$ws.Range("K32").Value = 1725
$ws.Range("L32").Value = 1701
$ws.Range("M32").Value = -1399
$ws.Range("N32").Value = -2353
$ws.Range("H41").Value = 5848382.5
$ws.Range("I41").Value = 15873130
$ws.Range("J41").Value = 612.8333
$ws.Range("K41").Value = 15873130
$ws.Range("L41").Value = 612.8333
$ws.Range("M41").Value = -15872690
$ws.Range("N41").Value = -1492.8333
$ws.Range("H46").Value = 1900
$ws.Range("J46").Value = 1900
$ws.Range("L46").Value = 5700
$ws.Range("N46").Value = -5938
$ws.Range("H55").Value = 509.76923
$ws.Range("I55").Value = 100
$ws.Range("J55").Value = 543.9167
$ws.Range("K55").Value = 100
$ws.Range("L55").Value = 543.9167
$ws.Range("M55").Value = 114
$ws.Range("N55").Value = -971.9167
$ws.Range("H60").Value = 1900
$ws.Range("J60").Value = 1900
$ws.Range("L60").Value = 5700
$ws.Range("N60").Value = -6668
$ws.Range("H95").Value = 25000
$ws.Range("J95").Value = 25000
$ws.Range("L95").Value = 25000
$ws.Range("N95").Value = -30492
$ws.Range("H107").Value = 757.875
$ws.Range("I107").Value = 848.95
$ws.Range("K107").Value = 848.95
$ws.Range("M107").Value = 1071.05
$ws.Range("H129").Value = 891
$ws.Range("I129").Value = 259.22223
$ws.Range("J129").Value = 1225.4706
$ws.Range("K129").Value = 777.66669
$ws.Range("L129").Value = 3676.4118
$ws.Range("M129").Value = 4222.33331
$ws.Range("N129").Value = -13676.4118
$ws.Range("H135").Value = 5544.1055
$ws.Range("I135").Value = 4096
$ws.Range("K135").Value = 36864
$ws.Range("M135").Value = -34329
$ws.Range("H141").Value = 5445.385
$ws.Range("I141").Value = 5574.1665
$ws.Range("J141").Value = 3900
$ws.Range("K141").Value = 16722.4995
$ws.Range("L141").Value = 11700
$ws.Range("M141").Value = -11542.4995
$ws.Range("N141").Value = -22060

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4207.2856
$ws.Range("I2").Value = 2987.889
$ws.Range("J2").Value = 6402.2
$ws.Range("K2").Value = 2987.889
$ws.Range("L2").Value = 6402.2
$ws.Range("M2").Value = -2874.889
$ws.Range("N2").Value = -6628.2
$ws.Range("H30").Value = 4950
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 4950
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 4950
$ws.Range("M30").ClearContents()
$ws.Range("N30").Value = -5250
$ws.Range("H32").Value = 363443.06
$ws.Range("I32").Value = 3456.2112
$ws.Range("J32").Value = 3203339.2
$ws.Range("K32").Value = 3456.2112
$ws.Range("L32").Value = 3203339.2
$ws.Range("M32").Value = -3169.2112
$ws.Range("N32").Value = -3203913.2
$ws.Range("H97").Value = 995.6
$ws.Range("I97").Value = 883.64514
$ws.Range("J97").Value = 1381.2222
$ws.Range("K97").Value = 883.64514
$ws.Range("L97").Value = 1381.2222
$ws.Range("M97").Value = -387.64514
$ws.Range("N97").Value = -2373.2222
$ws.Range("H110").Value = 776.5454999999999
$ws.Range("I110").Value = 701.4706
$ws.Range("J110").Value = 1031.8
$ws.Range("K110").Value = 701.4706
$ws.Range("L110").Value = 1031.8
$ws.Range("M110").Value = 1343.5294
$ws.Range("N110").Value = -5121.8
$ws.Range("H116").Value = 4207.2856
$ws.Range("I116").Value = 2987.889
$ws.Range("J116").Value = 6402.2
$ws.Range("K116").Value = 2987.889
$ws.Range("L116").Value = 6402.2
$ws.Range("M116").Value = -693.8890000000001
$ws.Range("N116").Value = -10990.2
$ws.Range("H122").Value = 41420.46
$ws.Range("I122").Value = 2837.6191
$ws.Range("K122").Value = 8512.8573
$ws.Range("M122").Value = -6062.8573
$ws.Range("H131").Value = 60558.43
$ws.Range("J131").Value = 60558.43
$ws.Range("L131").Value = 60558.43
$ws.Range("N131").Value = -70638.42999999999
$ws.Range("H132").Value = 17876384
$ws.Range("I132").Value = 25001186
$ws.Range("J132").Value = 64380.25
$ws.Range("K132").Value = 75003558
$ws.Range("L132").Value = 193140.75
$ws.Range("M132").Value = -75001028
$ws.Range("N132").Value = -198200.75

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4207.2856
$ws.Range("I3").Value = 2987.889
$ws.Range("J3").Value = 6402.2
$ws.Range("K3").Value = 2987.889
$ws.Range("L3").Value = 6402.2
$ws.Range("M3").Value = -2873.889
$ws.Range("N3").Value = -6630.2
$ws.Range("H107").Value = 936.9394
$ws.Range("I107").Value = 623.86957
$ws.Range("J107").Value = 1657
$ws.Range("K107").Value = 623.86957
$ws.Range("L107").Value = 1657
$ws.Range("M107").Value = 1296.13043
$ws.Range("N107").Value = -5497

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1286.5
$ws.Range("I16").Value = 1176.375
$ws.Range("J16").Value = 1433.3334
$ws.Range("K16").Value = 1176.375
$ws.Range("L16").Value = 1433.3334
$ws.Range("M16").Value = -889.375
$ws.Range("N16").Value = -2007.3334
$ws.Range("H31").Value = 4233.3716
$ws.Range("I31").Value = 6542.1577
$ws.Range("K31").Value = 6542.1577
$ws.Range("M31").Value = -6247.1577
$ws.Range("H34").Value = 4233.3716
$ws.Range("I34").Value = 6542.1577
$ws.Range("K34").Value = 6542.1577
$ws.Range("M34").Value = -6340.1577
$ws.Range("H58").Value = 2781.1428
$ws.Range("I58").Value = 1156
$ws.Range("K58").Value = 1156
$ws.Range("M58").Value = -953
$ws.Range("H105").Value = 1253.2222
$ws.Range("I105").Value = 918.4286
$ws.Range("J105").Value = 2425
$ws.Range("K105").Value = 918.4286
$ws.Range("L105").Value = 2425
$ws.Range("M105").Value = 828.5714
$ws.Range("N105").Value = -5919
$ws.Range("H113").Value = 1286.5
$ws.Range("I113").Value = 1176.375
$ws.Range("J113").Value = 1433.3334
$ws.Range("K113").Value = 1176.375
$ws.Range("L113").Value = 1433.3334
$ws.Range("M113").Value = 993.625
$ws.Range("N113").Value = -5773.3334
$ws.Range("H122").Value = 1458.8572
$ws.Range("I122").Value = 706
$ws.Range("J122").Value = 1760
$ws.Range("K122").Value = 2118
$ws.Range("L122").Value = 5280
$ws.Range("M122").Value = 332
$ws.Range("N122").Value = -10180
$ws.Range("H131").Value = 30000
$ws.Range("J131").Value = 30000
$ws.Range("L131").Value = 30000
$ws.Range("N131").Value = -40080
$ws.Range("H134").Value = 2465.9412
$ws.Range("I134").Value = 1541.4166
$ws.Range("J134").Value = 4684.8
$ws.Range("K134").Value = 4624.2498
$ws.Range("L134").Value = 14054.4
$ws.Range("M134").Value = -2089.2498
$ws.Range("N134").Value = -19124.4
$ws.Range("H136").Value = 2781.1428
$ws.Range("I136").Value = 1156
$ws.Range("K136").Value = 3468
$ws.Range("M136").Value = -918

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1090.4166
$ws.Range("I5").Value = 1090.4166
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3271.2498
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3159.2498
$ws.Range("N5").ClearContents()
$ws.Range("H14").Value = 71.5
$ws.Range("I14").Value = 71.5
$ws.Range("K14").Value = 214.5
$ws.Range("M14").Value = -41.5
$ws.Range("H109").Value = 1229.125
$ws.Range("I109").Value = 958.5
$ws.Range("J109").Value = 1499.75
$ws.Range("K109").Value = 2875.5
$ws.Range("L109").Value = 4499.25
$ws.Range("M109").Value = -1835.5
$ws.Range("N109").Value = -6579.25
$ws.Range("H132").Value = 1874.65
$ws.Range("I132").Value = 1120.6
$ws.Range("J132").Value = 2126
$ws.Range("K132").Value = 10085.4
$ws.Range("L132").Value = 19134
$ws.Range("M132").Value = -7555.4
$ws.Range("N132").Value = -24194
$ws.Range("H135").Value = 1090.4166
$ws.Range("I135").Value = 1090.4166
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9813.749400000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -7278.749400000001
$ws.Range("N135").ClearContents()

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 5002002.5
$ws.Range("J20").Value = 4005
$ws.Range("L20").Value = 4005
$ws.Range("N20").Value = -4495
$ws.Range("H102").Value = 1354.9395
$ws.Range("I102").Value = 1055.1923
$ws.Range("J102").Value = 2468.2856
$ws.Range("K102").Value = 1055.1923
$ws.Range("L102").Value = 2468.2856
$ws.Range("M102").Value = 566.8077000000001
$ws.Range("N102").Value = -5712.2856

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2995
$ws.Range("J20").Value = 2995
$ws.Range("L20").Value = 2995
$ws.Range("N20").Value = -3447
$ws.Range("H22").Value = 745.6429000000001
$ws.Range("I22").Value = 509.5238
$ws.Range("J22").Value = 981.7619
$ws.Range("K22").Value = 509.5238
$ws.Range("L22").Value = 981.7619
$ws.Range("M22").Value = -214.5238
$ws.Range("N22").Value = -1571.7619
$ws.Range("H27").Value = 745.6429000000001
$ws.Range("I27").Value = 509.5238
$ws.Range("J27").Value = 981.7619
$ws.Range("K27").Value = 509.5238
$ws.Range("L27").Value = 981.7619
$ws.Range("M27").Value = -402.5238
$ws.Range("N27").Value = -1195.7619
$ws.Range("H136").Value = 9277.333000000001
$ws.Range("I136").Value = 8883.5
$ws.Range("J136").Value = 10537.6
$ws.Range("K136").Value = 26650.5
$ws.Range("L136").Value = 31612.8
$ws.Range("M136").Value = -24100.5
$ws.Range("N136").Value = -36712.8

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").ClearContents()
$ws.Range("H8").Value = 2003
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").ClearContents()
$ws.Range("H81").Value = 3297.6875
$ws.Range("I81").Value = 2694.3
$ws.Range("J81").Value = 4303.3335
$ws.Range("K81").Value = 5388.6
$ws.Range("L81").Value = 8606.666999999999
$ws.Range("M81").Value = -4327.6
$ws.Range("N81").Value = -10728.667
$ws.Range("H84").Value = 3297.6875
$ws.Range("I84").Value = 2694.3
$ws.Range("J84").Value = 4303.3335
$ws.Range("K84").Value = 26943
$ws.Range("L84").Value = 43033.335
$ws.Range("M84").Value = -21639
$ws.Range("N84").Value = -53641.335
$ws.Range("H122").Value = 1303.6538
$ws.Range("I122").Value = 1347.3334
$ws.Range("J122").Value = 1244.091
$ws.Range("K122").Value = 4042.0002
$ws.Range("L122").Value = 3732.273
$ws.Range("M122").Value = -1592.0002
$ws.Range("N122").Value = -8632.272999999999
